$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: refresh the "as of" date in I1 (quote-prefixed so the COM layer
# keeps it as literal text instead of coercing it to a date serial)
$ws.Range("I1").Value = "'07/04/2023"

# Remove the "NOT" row (row 8) entirely - its category was an initialization
# placeholder that is no longer emitted; deleting the row shifts TEC/TST/VIP/WLC
# up by one and shrinks the used range from J12 to J11.
$ws.Rows.Item(8).Delete()

# Refresh the per-category counters (hour-on-activity snapshot updated to the
# new reporting hour). B column (pasthour flag) now only ever carries a value
# for hour 00, so every category resets to 0.
$rows = @(
    @("AMM",     0, 459.123, 674.1279999999999, 77, 8, 26, 136, 2139, -68.48396446937821),
    @("IPR",     0, 151,     154,                2,  1, 2,  0,   181,  -14.9171270718232),
    @("MIG",     0, 279,     285,                11, 0, 8,  0,   461,  -38.17787418655097),
    @("MOB",     0, 814,     934,                43, 6, 16, 71,  1302, -28.2642089093702),
    @("MOB PRE", 0, 469,     510,                41, 3, 9,  0,   1109, -54.01262398557258),
    @("MSK",     0, 227,     257,                30, 0, 3,  0,   289,  -11.07266435986159),
    @("TEC",     0, 402,     406,                3,  1, 3,  0,   905,  -55.13812154696133),
    @("TST",     0, 135,     143,                10, 0, 3,  0,   141,  1.418439716312059),
    @("VIP",     0, 4,       7,                  3,  0, 0,  0,   4,    75),
    @("WLC",     0, 31,      32,                 1,  0, 0,  0,   96,   -66.66666666666667)
)

$rowIndex = 2
foreach ($r in $rows) {
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $r[$col - 1]
    }
    $rowIndex = $rowIndex + 1
}
